# Adds the new quarterly row (01-07-2021) to the bottom of the data table,
# matching the "Actualización desde MV -datos-" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 76

# Column A holds a text label formatted like a date ("01-07-2021"), not an
# actual date value. Force the cell to Text first so the literal string is
# stored verbatim (as a shared string) instead of being auto-converted into
# a date serial number, then clear the format back to Normal so the cell
# ends up with no explicit style, same as the rest of the column.
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "01-07-2021"
$ws.Range("A$newRow").Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = 14232
$ws.Cells.Item($newRow, 3).Value = 2003
$ws.Cells.Item($newRow, 4).Value = 862
$ws.Cells.Item($newRow, 5).Value = 656
$ws.Cells.Item($newRow, 6).Value = 486
$ws.Cells.Item($newRow, 7).Value = 1520
$ws.Cells.Item($newRow, 8).Value = 2837
$ws.Cells.Item($newRow, 9).Value = -1317
$ws.Cells.Item($newRow, 10).Value = -2071
$ws.Cells.Item($newRow, 11).Value = 755
$ws.Cells.Item($newRow, 12).Value = -2371
$ws.Cells.Item($newRow, 13).Value = 4406
$ws.Cells.Item($newRow, 14).Value = -58
$ws.Cells.Item($newRow, 15).Value = -84
$ws.Cells.Item($newRow, 16).Value = 4573
$ws.Cells.Item($newRow, 17).Value = -25
$ws.Cells.Item($newRow, 18).Value = 8675
